$wb = $excel.ActiveWorkbook

# --- Sheet "2o Parcial" ---
$ws2 = $wb.Worksheets.Item("2o Parcial")

# Row 5 (Miguel Lopez Yadira / 2AEV)
$ws2.Range("E5").Value = 34
$ws2.Range("F5").Value = 1
$ws2.Range("G5").Value = 97.09999999999999
$ws2.Range("H5").Value = 2.9
$ws2.Range("I5").Value = 6.9
$ws2.Range("J5").Value = 0
$ws2.Range("K5").Value = 0

# Row 6 (Miguel Lopez Yadira / Cultura digital II)
$ws2.Range("E6").Value = 34
$ws2.Range("F6").Value = 1
$ws2.Range("G6").Value = 97.09999999999999
$ws2.Range("H6").Value = 2.9
$ws2.Range("I6").Value = 6.9
$ws2.Range("J6").Value = 0
$ws2.Range("K6").Value = 0

# Row 16 (Totales Generales)
$ws2.Range("E16").Value = 330
$ws2.Range("F16").Value = 29
$ws2.Range("G16").Value = 91.90000000000001
$ws2.Range("H16").Value = 8.1
$ws2.Range("I16").Value = 8.1
$ws2.Range("J16").Value = 0
$ws2.Range("K16").Value = 0

# --- Sheet "Final" ---
$wsFinal = $wb.Worksheets.Item("Final")

# Row 5 (Miguel Lopez Yadira / 2AEV) - Promedio updated
$wsFinal.Range("I5").Value = 7.3

# Row 6 (Miguel Lopez Yadira / Cultura digital II) - Promedio updated
$wsFinal.Range("I6").Value = 7.3
